# Add test on displaying detail by clicking name of tank on the list
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Location")

# Row 14 continues the alternating row banding - copy formatting from row 12
# (the previous "shaded" row) down to the new row 14, then fill in the values.
$ws.Range("A12:G12").Copy()
$ws.Range("A14:G14").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(14, 1).Value = "Detail.Name"
$ws.Cells.Item(14, 2).Value = "XPath"
$ws.Cells.Item(14, 3).Value = "//dt[text()='Name']/following-sibling::dd[1]"
$ws.Cells.Item(14, 4).Value = $null
$ws.Cells.Item(14, 5).Value = $null
$ws.Cells.Item(14, 6).Value = $null
$ws.Cells.Item(14, 7).Value = $null

$ws.Rows.Item(14).RowHeight = 20.25
